# tex-be atiras, es haladas
#
# The "A jelen peldaban ..." paragraph's text is extended with two more
# sentences (now spread across three runs instead of one), and the
# _GoBack bookmark that used to sit at the end of that paragraph moves to
# the empty paragraph right after "--todo" (the one right before the
# trailing blank paragraph / section break).

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphIndex($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    throw "Paragraph starting with '$needle' not found"
}

# --- Step 1: rewrite the "A jelen peldaban ..." paragraph's runs ----------

$run1 = "A jelen példában 3 részre van bontva a kotta, szöveg, akkord, pozíció. Ahhoz, hogy ebből tudjon  működni a konzolra kiiratás a megfelelő pozíciókkal, az indexelést tettem lista adatstrukt"
$run2 = "úrába, hogy meglegyen mindegyik akkordnak a megfelelő pozíciója. Kétfajta adat kerül tárolásra, az egyik mindenképpen az, hogy az xml-ben megadott string alapján melyik indexű karakternél található a szegmens szövegében a string maga, mert az jelöli az akkord helyét."
$run3 = " A másik pedig maga a string ami kiadja a kotta véglegesítését, a konkatenált akkord, - és szóközszámmal, valamint a kotta szövegével."

$idx1 = Find-ParagraphIndex("A jelen p")
$p1 = $d.Paragraphs.Item($idx1)

$xml1 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Standard`"/><w:ind w:left=`"360`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/><w:u w:val=`"single`"/></w:rPr></w:pPr>" +
        "<w:r><w:t>$run1</w:t></w:r>" +
        "<w:r><w:t>$run2</w:t></w:r>" +
        "<w:r><w:t xml:space=`"preserve`">$run3</w:t></w:r></w:p>"

$p1.Range.InsertXML($xml1) | Out-Null

# --- Step 2: move the _GoBack bookmark onto the paragraph right after -----
# --- "--todo" (was empty, no bookmark before) ------------------------------

$idxTodo = Find-ParagraphIndex("--todo")
$p2 = $d.Paragraphs.Item($idxTodo + 1)

$xml2 = "<w:p $wns><w:pPr><w:pStyle w:val=`"Standard`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr>" +
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$p2.Range.InsertXML($xml2) | Out-Null

Write-Host "Applied edit."
